$d = $word.ActiveDocument

# --- 1. "...equal to the content parameter." -> "...equal to the content property." (paragraph ~47)
$r = $d.Paragraphs(47).Range
$r.Find.Execute(" parameter", $true, $false, $false, $false, $false, $true, 1, $false, " property", 2)

# --- 2. "...which allows the id parameter to be accessed..." -> "...the id property to be accessed..." (paragraph ~55)
$r = $d.Paragraphs(55).Range
$r.Find.Execute("id parameter", $true, $false, $false, $false, $false, $true, 1, $false, "id property", 2)

# --- 3. "The edit view expects a model parameter called todos" -> "...model property called todos" (paragraph ~56)
$r = $d.Paragraphs(56).Range
$r.Find.Execute("model parameter called", $true, $false, $false, $false, $false, $true, 1, $false, "model property called", 2)

# --- 4. "The view also expects a model parameter called current" -> "...model property called current" (paragraph ~57)
#     (note: same paragraph also has "input parameter from the URI" which must stay unchanged)
$r = $d.Paragraphs(57).Range
$r.Find.Execute("expects a model parameter called", $true, $false, $false, $false, $false, $true, 1, $false, "expects a model property called", 2)

# --- 5. "The model map should have two parameters, todos" -> "...two properties, todos" (paragraph ~58)
$r = $d.Paragraphs(58).Range
$r.Find.Execute("two parameters", $true, $false, $false, $false, $false, $true, 1, $false, "two properties", 2)

# --- 6. "if the id parameter matches the current" -> "if the id property matches the current" (paragraph ~67)
#     (search excludes the leading "id" on purpose -- that "id" is in its own InlineCode-styled
#     run, and including it in the match would bleed the InlineCode style into the replacement)
$r = $d.Paragraphs(67).Range
$r.Find.Execute(" parameter matches the current", $true, $false, $false, $false, $false, $true, 1, $false, " property matches the current", 2)

# --- 7. "...with the content parameter, then..." -> "...with the content property, then..." (paragraph ~67)
$r = $d.Paragraphs(67).Range
$r.Find.Execute(" parameter, then", $true, $false, $false, $false, $false, $true, 1, $false, " property, then", 2)

# --- 8. Move the _GoBack bookmark from the end of paragraph 56 ("...passing the model map.")
#     to paragraph 71, right after "a named" (before " parameter for the id").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$p71 = $d.Paragraphs(71)
$anchor = $p71.Range.Duplicate
$anchor.Find.Execute("a named", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $target)
